$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rates text block in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$rA1 = $wsHoja1.Range("A1")
$oldText = $rA1.Value()
$newText = $oldText -replace [regex]::Escape("✅ 1000 Bs = 1.87 = 6934.14 pesos`n✅ 6934.14 pesos = 1.86 = 944.89 Bs"), "✅ 1000 Bs = 1.79 = 6607.32 pesos`n✅ 6607.32 pesos = 1.78 = 940.6 Bs"
$rA1.Value() = $newText

# --- Sheet "tasas": update N10, O10, N12, O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value() = 560
$wsTasas.Range("O10").Value() = 3700.1
$wsTasas.Range("N12").Value() = 3715.99
$wsTasas.Range("O12").Value() = 529
